$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the Area / Atotal columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area per segment: first row referenced against 0, rest against previous depth
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Total area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Match the author's final selection (cell H2)
$ws.Range("H2").Select() | Out-Null
